$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I header text (style copied from H4: bold, centered, wrapped)
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = "Penalty Issued"

# New column I totals-row cell (style copied from H11: bold, bottom border, centered)
$ws.Range("H11").Copy()
$ws.Range("I11").PasteSpecial(-4122)

# New column I data placeholder row (style copied from H6, no style -> plain)
$ws.Range("I6").Value = "{d.Reg[i].Penalty}"

# Thin bottom-border-only cell in the spacer row (I3), then reuse the same
# resulting style for I5 so no duplicate style entries are produced.
$ws.Range("I3").Borders.Item(9).Color = 0
$ws.Range("I3").Borders.Item(9).ColorIndex = 64
$ws.Range("I3").Borders.Item(9).LineStyle = 1
$ws.Range("I3").Copy()
$ws.Range("I5").PasteSpecial(-4122)

# Update the active selection to match the saved view state
$ws.Range("F10").Select()
